# Apply edit described by the diff:
# - D3 and D4 change text from "always caps" to "change word"
# - D3 and D4 adopt the same style/format as column C (s=2) instead of s=1
# - A new shared string "change word" is introduced (Excel will manage the
#   shared string table reordering automatically when saving)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from C3/C4 (style index 2) onto D3/D4 so the cells end
# up using the same cell style as the rest of that column, matching the
# target workbook.
$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C4").Copy()
$ws.Range("D4").PasteSpecial(-4122)  # xlPasteFormats

# Update the cell text/value for D3 and D4.
$ws.Range("D3").Value = "change word"
$ws.Range("D4").Value = "change word"

$excel.CutCopyMode = $false
